# Add a new weekly survival-count column (AA) that duplicates the most
# recent week's counts (column Z) for every data row, and move the active
# selection to reflect the new working cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 37; $row++) {
    $zValue = $ws.Cells.Item($row, 26).Value2   # column Z = 26
    $ws.Cells.Item($row, 27).Value = $zValue    # column AA = 27
}

$ws.Range("AC7").Select()
